# TAB-8861 - Support online event fields in SGT spreadsheet imports
#
# Adds three new columns (Delivery Method / Online Delivery URL / Online
# Platform) to the "Events" sheet of the SGT import template, populates a
# handful of sample rows with the new fields (including two hyperlinked
# "Online Delivery URL" cells), and moves the sheet's view/selection to
# reflect where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# --- New header cells -------------------------------------------------
$ws.Range("N1").Value = "Delivery Method"
$ws.Range("O1").Value = "Online Delivery URL"
$ws.Range("P1").Value = "Online Platform"

# --- Sample data for the new columns -----------------------------------
$ws.Range("N3").Value = "Hybrid"
$ws.Range("O3").Value = "https://tabula.warwick.ac.uk"
$ws.Range("P3").Value = "Teams"

# Row 4 previously held a (now stale) Location value in J4 - the import
# moves that row to be an online-only event, so the old J4 value is
# cleared in favour of the new Delivery Method/URL/Platform fields.
$ws.Range("J4").ClearContents()
$ws.Range("N4").Value = "OnlineOnly"
$ws.Range("O4").Value = "https://tabula.warwick.ac.uk"
$ws.Range("P4").Value = "Moodle"

$ws.Range("N5").Value = "FaceToFaceOnly"

$ws.Range("N8").Value = "Junk"
$ws.Range("P9").Value = "Junk"

# --- Hyperlinks for the new "Online Delivery URL" cells ----------------
$ws.Hyperlinks.Add($ws.Range("O3"), "https://tabula.warwick.ac.uk", "", "", "https://tabula.warwick.ac.uk")
$ws.Hyperlinks.Add($ws.Range("O4"), "https://tabula.warwick.ac.uk", "", "", "https://tabula.warwick.ac.uk")

# --- Column widths for the new columns ----------------------------------
$ws.Columns("M").ColumnWidth = 10.95
$ws.Columns("N").ColumnWidth = 16.57
$ws.Columns("O").ColumnWidth = 24.84
$ws.Columns("P").ColumnWidth = 19.77
$ws.Columns("Q").ColumnWidth = 20.21

# --- View state: scroll across to show the new columns and select the
# newly-populated hyperlink cell, as the importer author last left it.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("O3").Select()
